# Update building names: replace comma separator with a dash/hyphen
# for the two "SOUTHBANK PERFORMING ARTS" rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "SOUTHBANK PERFORMING ARTS - ST KILDA RD"
$ws.Range("A16").Value = "SOUTHBANK PERFORMING ARTS - DODDS ST"

# Match the resulting active selection in the saved file (A16)
$ws.Range("A16").Select()
